$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate the existing "High Priority break-up" sheet (sheet index 5) ---
# The copy preserves all formatting/styles and becomes the new last sheet.
$wsOld = $wb.Worksheets.Item(5)
$wsOld.Copy([System.Reflection.Missing]::Value, $wsOld)

# The newly created copy is placed right after $wsOld, i.e. at index 6.
$wsNew = $wb.Worksheets.Item(6)
$wsNew.Name = "Major update - High Priority "

# --- Step 2: rename the original sheet and replace its data with the new figures ---
$wsOld.Name = "Interannual update - High Pri"

$wsOld.Cells.Item(2,1).Value = "Trend New"
$wsOld.Cells.Item(2,2).Value = 88
$wsOld.Cells.Item(2,3).Value = 85.40000000000001
$wsOld.Cells.Item(2,4).Value = 88
$wsOld.Cells.Item(2,5).Value = 98.90000000000001

$wsOld.Cells.Item(3,1).Value = "Trend Different"
$wsOld.Cells.Item(3,2).Value = 4
$wsOld.Cells.Item(3,3).Value = 3.9
$wsOld.Cells.Item(3,4).ClearContents()
$wsOld.Cells.Item(3,5).ClearContents()

$wsOld.Cells.Item(4,1).Value = "IUCN"
$wsOld.Cells.Item(4,2).Value = 11
$wsOld.Cells.Item(4,3).Value = 10.7
$wsOld.Cells.Item(4,4).Value = 1
$wsOld.Cells.Item(4,5).Value = 1.1

# --- Step 3: update "Trends Status" sheet (index 1) ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2,2).Value = 5
$ws1.Cells.Item(2,3).Value = 2
$ws1.Cells.Item(2,4).Value = 41.7
$ws1.Cells.Item(2,5).Value = 9.1

$ws1.Cells.Item(3,2).Value = 4
$ws1.Cells.Item(3,3).Value = 9
$ws1.Cells.Item(3,4).Value = 33.3
$ws1.Cells.Item(3,5).Value = 40.9

$ws1.Cells.Item(4,2).Value = 2
$ws1.Cells.Item(4,3).Value = 9
$ws1.Cells.Item(4,4).Value = 16.7
$ws1.Cells.Item(4,5).Value = 40.9

$ws1.Cells.Item(5,3).Value = 1
$ws1.Cells.Item(5,5).Value = 4.5

$ws1.Cells.Item(6,2).Value = 1
$ws1.Cells.Item(6,3).Value = 1
$ws1.Cells.Item(6,4).Value = 8.300000000000001
$ws1.Cells.Item(6,5).Value = 4.5

$ws1.Cells.Item(7,2).Value = 19
$ws1.Cells.Item(7,3).Value = 35

# --- Step 4: update "Priority Status" sheet (index 3) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,2).Value = 103
$ws3.Cells.Item(3,2).Value = 286
$ws3.Cells.Item(4,2).Value = 554

# --- Step 5: update "Species qualification" sheet (index 4) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2,1).Value = "SoIB Assessment"
$ws4.Cells.Item(3,3).Value = 12
$ws4.Cells.Item(4,3).Value = 22
